# Updates the cryptos sheet with refreshed prices / % volume figures,
# and reorders three coin-pairs that swapped rank (rows 18/19, 24/25, 49/50).
#
# Commit: "Updated cryptos list on Wed Sep 11 13:51:11 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin names, links and % volume figures (and the handful of Price cells that
# aren't plain decimals, e.g. "56.256.52" or "0.0₃0714") can be assigned
# directly -- Excel's COM Value setter stores them as text on its own.
$plainUpdates = @{
    'D2' = '56.256.52'
    'E2' = '  -0.27%  '
    'D3' = '2.312.95'
    'E3' = '  -0.13%  '
    'E4' = '  -0.01%  '
    'E5' = '  -0.55%  '
    'E6' = '  -2.31%  '
    'E7' = '  +0.45%  '
    'E8' = '  -1.04%  '
    'E9' = '  -2.65%  '
    'E10' = '  +0.09%  '
    'E11' = '  -1.09%  '
    'E12' = '  -1.84%  '
    'E13' = '  -1.60%  '
    'D14' = '2.726.46'
    'E14' = '  -0.28%  '
    'D15' = '56.287.35'
    'E15' = '  -0.43%  '
    'E16' = '  -1.56%  '
    'D17' = '2.319.51'
    'E17' = '  -0.75%  '
    'B18' = 'BitcoinCash'
    'C18' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'E18' = '  +2.11%  '
    'B19' = 'Chainlink'
    'C19' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E19' = '  -0.79%  '
    'E20' = '  -1.86%  '
    'E21' = '  +1.61%  '
    'E22' = '  -0.03%  '
    'E23' = '  +0.19%  '
    'B24' = 'InternetComputer(DFINITY)'
    'C24' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'E24' = '  +8.54%  '
    'B25' = 'Kaspa'
    'C25' = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    'E25' = '  +0.60%  '
    'E26' = '  +0.41%  '
    'E27' = '  +1.85%  '
    'E28' = '  +1.21%  '
    'E29' = '  -0.95%  '
    'D30' = '0.0₃0714'
    'E30' = '  -3.13%  '
    'E31' = '  -0.97%  '
    'E33' = '  -0.01%  '
    'E34' = '  +0.54%  '
    'E35' = '  -1.05%  '
    'E36' = '  -2.57%  '
    'E37' = '  -4.49%  '
    'E38' = '  +1.28%  '
    'E39' = '  +1.84%  '
    'E40' = '  +6.38%  '
    'E41' = '  -2.29%  '
    'E42' = '  -0.48%  '
    'E43' = '  -0.02%  '
    'E44' = '  -4.53%  '
    'E45' = '  -0.70%  '
    'E46' = '  -2.39%  '
    'E47' = '  -1.86%  '
    'E48' = '  +1.74%  '
    'B49' = 'Polygon'
    'C49' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E49' = '  -0.07%  '
    'B50' = 'VeChain'
    'C50' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E50' = '  -1.26%  '
    'E51' = '  +1.38%  '
}

foreach ($cellRef in $plainUpdates.Keys) {
    $ws.Range($cellRef).Value = $plainUpdates[$cellRef]
}

# Price cells whose new value is a plain decimal number (e.g. "1.00", "6.67")
# must stay TEXT, matching the source sheet (every Price cell is stored as a
# string). A bare assignment would let Excel auto-coerce these into Number
# cells (and drop significant trailing zeros, e.g. "1.00" -> 1). Force text
# via NumberFormat "@", then restore the cell to the default "Normal" style
# so no stray number-format override is left behind.
$textUpdates = @{
    'D5' = '515.15'
    'D6' = '131.02'
    'D9' = '0.0998'
    'D11' = '5.22'
    'D12' = '0.335'
    'D13' = '23.42'
    'D18' = '329.61'
    'D19' = '10.35'
    'D21' = '6.67'
    'D22' = '1.00'
    'D23' = '60.85'
    'D24' = '8.62'
    'D25' = '0.164'
    'D26' = '0.998'
    'D28' = '168.30'
    'D31' = '6.11'
    'D32' = '18.27'
    'D34' = '0.999'
    'D36' = '3.91'
    'D37' = '0.879'
    'D39' = '38.59'
    'D40' = '148.26'
    'D43' = '279.91'
    'D45' = '0.0925'
    'D47' = '0.552'
    'D48' = '18.10'
    'D49' = '0.380'
    'D50' = '0.0214'
    'D51' = '17.03'
}

foreach ($cellRef in $textUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$cellRef]
    $cell.Style = "Normal"
}
